$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.35157
$ws.Range("H2").Value = 25.05471
$ws.Range("I2").Value = 0.3629556103554933
$ws.Range("J2").Value = 0.3629556103554933
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 14.80997966666667
$ws.Range("N2").Value = 44.429939
$ws.Range("O2").Value = 0.2388798507865045
$ws.Range("P2").Value = 0.2388798507865045
$ws.Range("Q2").Value = 123.6865818847434
$ws.Range("R2").Value = 1113.17923696269
$ws.Range("S2").Value = 0.08670278204384492
$ws.Range("T2").Value = 0.08670278204384493

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.35157
$ws.Range("H3").Value = 25.05471
$ws.Range("I3").Value = 0.3629556103554933
$ws.Range("J3").Value = 0.3629556103554933
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 33.42883866666667
$ws.Range("N3").Value = 100.286516
$ws.Range("O3").Value = 0.5391956081231261
$ws.Range("P3").Value = 0.5391956081231262
$ws.Range("Q3").Value = 279.1832861433734
$ws.Range("R3").Value = 2512.64957529036
$ws.Range("S3").Value = 0.1957040710473306
$ws.Range("T3").Value = 0.1957040710473307

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.35157
$ws.Range("H4").Value = 25.05471
$ws.Range("I4").Value = 0.3629556103554933
$ws.Range("J4").Value = 0.3629556103554933
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.509909333333333
$ws.Range("N4").Value = 7.529728
$ws.Range("O4").Value = 0.04048396962919451
$ws.Range("P4").Value = 0.04048396962919452
$ws.Range("Q4").Value = 20.96168349098667
$ws.Range("R4").Value = 188.65515141888
$ws.Range("S4").Value = 0.01469388390637755
$ws.Range("T4").Value = 0.01469388390637755

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.35157
$ws.Range("H5").Value = 25.05471
$ws.Range("I5").Value = 0.3629556103554933
$ws.Range("J5").Value = 0.3629556103554933
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 11.24888166666667
$ws.Range("N5").Value = 33.746645
$ws.Range("O5").Value = 0.1814405714611748
$ws.Range("P5").Value = 0.1814405714611748
$ws.Range("Q5").Value = 93.94582266088335
$ws.Range("R5").Value = 845.51240394795
$ws.Range("S5").Value = 0.06585487335794021
$ws.Range("T5").Value = 0.06585487335794023

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.216696
$ws.Range("H6").Value = 33.650088
$ws.Range("I6").Value = 0.4874727437897329
$ws.Range("J6").Value = 0.487472743789733
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 14.80997966666667
$ws.Range("N6").Value = 44.429939
$ws.Range("O6").Value = 0.2388798507865045
$ws.Range("P6").Value = 0.2388798507865045
$ws.Range("Q6").Value = 166.1190396871813
$ws.Range("R6").Value = 1495.071357184632
$ws.Range("S6").Value = 0.1164474162989793
$ws.Range("T6").Value = 0.1164474162989794

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.216696
$ws.Range("H7").Value = 33.650088
$ws.Range("I7").Value = 0.4874727437897329
$ws.Range("J7").Value = 0.487472743789733
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 33.42883866666667
$ws.Range("N7").Value = 100.286516
$ws.Range("O7").Value = 0.5391956081231261
$ws.Range("P7").Value = 0.5391956081231262
$ws.Range("Q7").Value = 374.9611209570454
$ws.Range("R7").Value = 3374.650088613408
$ws.Range("S7").Value = 0.2628431625311539
$ws.Range("T7").Value = 0.262843162531154

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.216696
$ws.Range("H8").Value = 33.650088
$ws.Range("I8").Value = 0.4874727437897329
$ws.Range("J8").Value = 0.487472743789733
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.509909333333333
$ws.Range("N8").Value = 7.529728
$ws.Range("O8").Value = 0.04048396962919451
$ws.Range("P8").Value = 0.04048396962919452
$ws.Range("Q8").Value = 28.15288997956266
$ws.Range("R8").Value = 253.376009816064
$ws.Range("S8").Value = 0.01973483175464366
$ws.Range("T8").Value = 0.01973483175464367

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.216696
$ws.Range("H9").Value = 33.650088
$ws.Range("I9").Value = 0.4874727437897329
$ws.Range("J9").Value = 0.487472743789733
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 11.24888166666667
$ws.Range("N9").Value = 33.746645
$ws.Range("O9").Value = 0.1814405714611748
$ws.Range("P9").Value = 0.1814405714611748
$ws.Range("Q9").Value = 126.1752859949733
$ws.Range("R9").Value = 1135.57757395476
$ws.Range("S9").Value = 0.088447333204956
$ws.Range("T9").Value = 0.08844733320495603

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.441627666666667
$ws.Range("H10").Value = 10.324883
$ws.Range("I10").Value = 0.1495716458547737
$ws.Range("J10").Value = 0.1495716458547737
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.80997966666667
$ws.Range("N10").Value = 44.429939
$ws.Range("O10").Value = 0.2388798507865045
$ws.Range("P10").Value = 0.2388798507865045
$ws.Range("Q10").Value = 50.97043576357078
$ws.Range("R10").Value = 458.733921872137
$ws.Range("S10").Value = 0.03572965244368024
$ws.Range("T10").Value = 0.03572965244368024

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.441627666666667
$ws.Range("H11").Value = 10.324883
$ws.Range("I11").Value = 0.1495716458547737
$ws.Range("J11").Value = 0.1495716458547737
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 33.42883866666667
$ws.Range("N11").Value = 100.286516
$ws.Range("O11").Value = 0.5391956081231261
$ws.Range("P11").Value = 0.5391956081231262
$ws.Range("Q11").Value = 115.0496160197365
$ws.Range("R11").Value = 1035.446544177628
$ws.Range("S11").Value = 0.08064837454464155
$ws.Range("T11").Value = 0.08064837454464156

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.441627666666667
$ws.Range("H12").Value = 10.324883
$ws.Range("I12").Value = 0.1495716458547737
$ws.Range("J12").Value = 0.1495716458547737
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 2.509909333333333
$ws.Range("N12").Value = 7.529728
$ws.Range("O12").Value = 0.04048396962919451
$ws.Range("P12").Value = 0.04048396962919452
$ws.Range("Q12").Value = 8.63817340242489
$ws.Range("R12").Value = 77.743560621824
$ws.Range("S12").Value = 0.006055253968173295
$ws.Range("T12").Value = 0.006055253968173297

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.441627666666667
$ws.Range("H13").Value = 10.324883
$ws.Range("I13").Value = 0.1495716458547737
$ws.Range("J13").Value = 0.1495716458547737
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 11.24888166666667
$ws.Range("N13").Value = 33.746645
$ws.Range("O13").Value = 0.1814405714611748
$ws.Range("P13").Value = 0.1814405714611748
$ws.Range("Q13").Value = 38.71446236305945
$ws.Range("R13").Value = 348.430161267535
$ws.Range("S13").Value = 0.0271383648982786
$ws.Range("T13").Value = 0.0271383648982786
